# Reverse the order of the comma-separated "Recorded By" values in column G.
# e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
# Only cells containing multiple comma-separated values are affected; single
# value cells are unchanged (reversing a single item is a no-op).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ",\s*"
        if ($parts.Count -gt 1) {
            $reversed = $parts[($parts.Count - 1)..0]
            $newVal = [string]::Join(", ", $reversed)
            $cell.Value2 = $newVal
        }
    }
}
